$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap data (columns A, Q, R, AC) between rows 3 and 4.
$a3 = $ws.Range("A3").Value2
$q3 = $ws.Range("Q3").Value2
$r3 = $ws.Range("R3").Value2
$ac3 = $ws.Range("AC3").Value2

$a4 = $ws.Range("A4").Value2
$q4 = $ws.Range("Q4").Value2
$r4 = $ws.Range("R4").Value2
$ac4 = $ws.Range("AC4").Value2

$ws.Range("A3").Value2 = $a4
$ws.Range("Q3").Value2 = $q4
$ws.Range("R3").Value2 = $r4
$ws.Range("AC3").Value2 = $ac4

$ws.Range("A4").Value2 = $a3
$ws.Range("Q4").Value2 = $q3
$ws.Range("R4").Value2 = $r3
$ws.Range("AC4").Value2 = $ac3

# Rotate data (columns A, Q, R, AC) among rows 13, 14, 15:
# new13 = old14, new14 = old15, new15 = old13
$a13 = $ws.Range("A13").Value2
$q13 = $ws.Range("Q13").Value2
$r13 = $ws.Range("R13").Value2
$ac13 = $ws.Range("AC13").Value2

$a14 = $ws.Range("A14").Value2
$q14 = $ws.Range("Q14").Value2
$r14 = $ws.Range("R14").Value2
$ac14 = $ws.Range("AC14").Value2

$a15 = $ws.Range("A15").Value2
$q15 = $ws.Range("Q15").Value2
$r15 = $ws.Range("R15").Value2
$ac15 = $ws.Range("AC15").Value2

$ws.Range("A13").Value2 = $a14
$ws.Range("Q13").Value2 = $q14
$ws.Range("R13").Value2 = $r14
$ws.Range("AC13").Value2 = $ac14

$ws.Range("A14").Value2 = $a15
$ws.Range("Q14").Value2 = $q15
$ws.Range("R14").Value2 = $r15
$ws.Range("AC14").Value2 = $ac15

$ws.Range("A15").Value2 = $a13
$ws.Range("Q15").Value2 = $q13
$ws.Range("R15").Value2 = $r13
$ws.Range("AC15").Value2 = $ac13
